# Auto-generated edit script: updates crypto price/volume table cells
# to match the target snapshot (GitHub Actions cryptos list update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'59.655.45"
$ws.Range("E2").Value = "  +0.49%  "
# Row 3
$ws.Range("D3").Value = "'2.650.53"
$ws.Range("E3").Value = "  +1.63%  "
# Row 4
$ws.Range("E4").Value = "  -0.06%  "
# Row 5
$ws.Range("D5").Value = "'537.50"
$ws.Range("E5").Value = "  -1.59%  "
# Row 6
$ws.Range("D6").Value = "'146.48"
$ws.Range("E6").Value = "  +3.69%  "
# Row 7
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.08%  "
# Row 8
$ws.Range("E8").Value = "  +1.21%  "
# Row 9
$ws.Range("D9").Value = "'6.78"
$ws.Range("E9").Value = "  +4.81%  "
# Row 10
$ws.Range("E10").Value = "  -0.45%  "
# Row 11
$ws.Range("E11").Value = "  +1.14%  "
# Row 12
$ws.Range("E12").Value = "  +0.08%  "
# Row 13
$ws.Range("D13").Value = "'3.118.66"
$ws.Range("E13").Value = "  +1.62%  "
# Row 14
$ws.Range("D14").Value = "'59.576.28"
$ws.Range("E14").Value = "  +0.50%  "
# Row 15
$ws.Range("E15").Value = "  +4.01%  "
# Row 16
$ws.Range("D16").Value = "'2.658.32"
$ws.Range("E16").Value = "  +2.35%  "
# Row 17
$ws.Range("E17").Value = "  +0.98%  "
# Row 18
$ws.Range("E18").Value = "  +2.51%  "
# Row 19
$ws.Range("D19").Value = "'340.05"
$ws.Range("E19").Value = "  -1.06%  "
# Row 20
$ws.Range("E20").Value = "  +2.27%  "
# Row 21
$ws.Range("D21").Value = "'6.19"
$ws.Range("E21").Value = "  -3.59%  "
# Row 22
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.02%  "
# Row 23
$ws.Range("D23").Value = "'66.58"
$ws.Range("E23").Value = "  -1.31%  "
# Row 24
$ws.Range("E24").Value = "  +2.37%  "
# Row 25
$ws.Range("E25").Value = "  -0.58%  "
# Row 26
$ws.Range("E26").Value = "  -0.01%  "
# Row 27
$ws.Range("D27").Value = "'7.30"
$ws.Range("E27").Value = "  +1.10%  "
# Row 28
$ws.Range("D28").Value = "'0.0₃0749"
$ws.Range("E28").Value = "  +1.46%  "
# Row 29
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.06%  "
# Row 30
$ws.Range("D30").Value = "'1.65"
$ws.Range("E30").Value = "  -3.21%  "
# Row 31
$ws.Range("D31").Value = "'5.87"
# Row 32
$ws.Range("D32").Value = "'18.91"
$ws.Range("E32").Value = "  +0.67%  "
# Row 33
$ws.Range("E33").Value = "  +1.13%  "
# Row 34
$ws.Range("E34").Value = "  +0.76%  "
# Row 35
$ws.Range("E35").Value = "  +2.54%  "
# Row 36
$ws.Range("E36").Value = "  +3.51%  "
# Row 37
$ws.Range("E37").Value = "  +0.24%  "
# Row 38
$ws.Range("E38").Value = "  -0.93%  "
# Row 39
$ws.Range("E39").Value = "  +1.74%  "
# Row 40
$ws.Range("D40").Value = "'286.79"
$ws.Range("E40").Value = "  +3.17%  "
# Row 41
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = "  -0.12%  "
# Row 42
$ws.Range("D42").Value = "'0.607"
$ws.Range("E42").Value = "  +1.52%  "
# Row 43
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'10.75"
$ws.Range("E43").Value = "  +0.28%  "
# Row 44
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "'0.0540"
$ws.Range("E44").Value = "  +3.04%  "
# Row 45
$ws.Range("D45").Value = "'19.31"
$ws.Range("E45").Value = "  +3.52%  "
# Row 46
$ws.Range("E46").Value = "  -0.99%  "
# Row 47
$ws.Range("E47").Value = "  +2.22%  "
# Row 48
$ws.Range("D48").Value = "'1.967.67"
$ws.Range("E48").Value = "  +1.01%  "
# Row 49
$ws.Range("E49").Value = "  +1.01%  "
# Row 50
$ws.Range("D50").Value = "'18.43"
$ws.Range("E50").Value = "  +0.46%  "
# Row 51
$ws.Range("D51").Value = "'111.33"
$ws.Range("E51").Value = "  +0.25%  "
